# "Cambios a noviembre 20 2023" — refresh Enfermeria month-over-month stats:
#  - Mes column: Sep-2023 (45170) -> Oct-2023 (45200) for every data row
#  - Procedimiento / Cantidad / Sede values updated per the new report
#  - A previously-unused duplicate string ("TOMA DE TENSÓN ARTERIAL", a typo
#    of "TOMA DE TENSIÓN ARTERIAL") drops out once nothing references it, and
#    a new procedure "LAVADO NASAL" appears (Bulevar, row 8)
#  - Selection cursor moves to C22

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# r, Sede, Mes(serial), Procedimiento, Cantidad
$rows = @(
    @(2,  "Bulevar",    45200, "INYECTOLOGÍA",             111),
    @(3,  "Bulevar",    45200, "TOMA DE EKG",               139),
    @(4,  "Bulevar",    45200, "LAVADO DE OÍDOS",            10),
    @(5,  "Bulevar",    45200, "RETIRO DE PUNTOS",            5),
    @(6,  "Bulevar",    45200, "TOMA DE TENSIÓN ARTERIAL",    9),
    @(7,  "Bulevar",    45200, "GLUCOMETRÍA",                 5),
    @(8,  "Bulevar",    45200, "LAVADO NASAL",                1),
    @(9,  "Bulevar",    45200, "CURACIÓN",                    3),
    @(10, "San Martin", 45200, "INYECTOLOGÍA",               73),
    @(11, "San Martin", 45200, "LAVADO DE OÍDOS",             6),
    @(12, "San Martin", 45200, "TOMA DE EKG ",              145),
    @(13, "San Martin", 45200, "RETIRO DE PUNTOS",            1),
    @(14, "San Martin", 45200, "CURACIÓN",                    9),
    @(15, "Cartagena",  45200, "CURACIÓN",                    6),
    @(16, "Cartagena",  45200, "GLUCOMETRÍA",                 1),
    @(17, "Cartagena",  45200, "INYECTOLOGÍA",               76),
    @(18, "Cartagena",  45200, "TOMA DE EKG",                94),
    @(19, "Cartagena",  45200, "TOMA DE TENSIÓN ARTERIAL",   95)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$ws.Range("C22").Select()

Write-Output "Enfermeria refreshed through October 2023"
